$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as text so values like "1.034"
# are not auto-converted to numbers by Excel, then restore default style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "19.930.16"
$ws.Range("E2").Value = "  -2.97%  "
$ws.Range("D3").Value = "1.413.19"
$ws.Range("E3").Value = "  -2.63%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.50%  "
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "276.29"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").Value = "0.3692"
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("D8").Value = "0.3105"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "39.90"
$ws.Range("D10").Value = "1.034"
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("D11").Value = "0.06512"
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").Value = "5.466"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").Value = "17.59"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").Value = "6.191"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "1.416.91"
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("D17").Value = "0.00001018"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "0.05689"
$ws.Range("E18").Value = "  -6.36%  "
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").Value = "70.88"
$ws.Range("E20").Value = "  -9.00%  "
$ws.Range("D21").Value = "5.594"
$ws.Range("E21").Value = "  -2.75%  "
$ws.Range("D22").Value = "14.72"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").Value = "10.96"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").Value = "2.233"
$ws.Range("E24").Value = "  -3.58%  "
$ws.Range("D25").Value = "19.967.16"
$ws.Range("E25").Value = "  -2.81%  "
$ws.Range("D26").Value = "2.265"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "133.02"
$ws.Range("E27").Value = "  -6.95%  "
$ws.Range("D28").Value = "17.18"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").Value = "1.577.64"
$ws.Range("E29").Value = "  -2.54%  "
$ws.Range("D30").Value = "110.04"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").Value = "3.920"
$ws.Range("E31").Value = "  +8.12%  "
$ws.Range("D32").Value = "5.199"
$ws.Range("E32").Value = "  -5.33%  "
$ws.Range("D33").Value = "0.8083"
$ws.Range("E33").Value = "  -12.51%  "
$ws.Range("D34").Value = "0.07769"
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("E35").Value = "  +2.75%  "
$ws.Range("D36").Value = "8.133"
$ws.Range("E36").Value = "  -1.95%  "
$ws.Range("D37").Value = "4.880"
$ws.Range("E37").Value = "  +1.58%  "
$ws.Range("D38").Value = "0.05832"
$ws.Range("E38").Value = "  +3.18%  "
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("D40").Value = "0.02042"
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("D41").Value = "10.42"
$ws.Range("E41").Value = "  -5.77%  "
$ws.Range("D42").Value = "1.105"
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("D43").Value = "0.1883"
$ws.Range("E43").Value = "  -2.19%  "
$ws.Range("D44").Value = "12.37"
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("D45").Value = "0.5300"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").Value = "3.534"
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("D47").Value = "116.44"
$ws.Range("E47").Value = "  +4.99%  "
$ws.Range("D48").Value = "0.5176"
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("D49").Value = "1.765"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("D50").Value = "1.033"
$ws.Range("E50").Value = "  -3.28%  "
$ws.Range("E51").Value = "  -0.49%  "

$ws.Range("D2:D51").Style = "Normal"
